# B2 R3 leachate weights and associated data
$wb = $excel.ActiveWorkbook

# The third tab ("Sheet1") is the active sheet - rename it to "Rainfall 3"
$ws = $wb.ActiveSheet
$ws.Name = "Rainfall 3"

# Select F20 as the active cell on this sheet (matches saved selection in diff)
$ws.Range("F20").Select()

# Data to fill in: rows 2-19, columns C (rain_date), F (water_plus_bottle_mass_collected),
# G (water_mass, formula F-E), H (DNA_filter_date)
$rainDate = 20210614

$fValues = @{
    2 = 906
    3 = 782
    4 = 810
    5 = 880
    6 = 794
    7 = 833
    8 = 882
    9 = 788
    10 = 829
    11 = 850
    12 = 915
    13 = 854
    14 = 840
    15 = 852
    16 = 887
    17 = 884
    18 = 895
    19 = 771
}

for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 3).Value = $rainDate
    $ws.Cells.Item($row, 6).Value = $fValues[$row]
    $ws.Cells.Item($row, 8).Value = $rainDate
}

# G2 gets its own (non-shared) formula; G3:G19 becomes a shared formula block.
$ws.Range("G2").Formula = "=F2-E2"
$ws.Range("G3:G19").Formula = "=F3-E3"
